# Update team-specific matrix values (Eastern Wash._A) from games pulled March 7.
# Applies updated probability figures to the existing matrix on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1542857142857143
$ws.Range("C2").Value = 0.64
$ws.Range("J2").Value = 0.01714285714285714
$ws.Range("P2").Value = 0.12
$ws.Range("S2").Value = 0.06857142857142857
$ws.Range("B3").Value = 0.004219409282700422
$ws.Range("C3").Value = 0.0379746835443038
$ws.Range("J3").Value = 0.01687763713080169
$ws.Range("P3").Value = 0.7721518987341772
$ws.Range("S3").Value = 0.1687763713080169
$ws.Range("J4").Value = 0.07407407407407407
$ws.Range("P4").Value = 0.6481481481481481
$ws.Range("S4").Value = 0.2777777777777778
$ws.Range("B6").Value = 0.09
$ws.Range("D6").Value = 0.01
$ws.Range("F6").Value = 0.09333333333333334
$ws.Range("J6").Value = 0.22
$ws.Range("O6").Value = 0.04
$ws.Range("Q6").Value = 0.1733333333333333
$ws.Range("R6").Value = 0.05333333333333334
$ws.Range("S6").Value = 0.32
$ws.Range("B7").Value = 0.1111111111111111
$ws.Range("D7").Value = 0.004629629629629629
$ws.Range("F7").Value = 0.07407407407407407
$ws.Range("J7").Value = 0.1388888888888889
$ws.Range("O7").Value = 0.02777777777777778
$ws.Range("Q7").Value = 0.1759259259259259
$ws.Range("R7").Value = 0.09722222222222222
$ws.Range("S7").Value = 0.3703703703703703
$ws.Range("B8").Value = 0.09791666666666667
$ws.Range("D8").Value = 0.03125
$ws.Range("F8").Value = 0.07708333333333334
$ws.Range("J8").Value = 0.075
$ws.Range("O8").Value = 0.03541666666666667
$ws.Range("Q8").Value = 0.1895833333333333
$ws.Range("R8").Value = 0.1125
$ws.Range("S8").Value = 0.38125
$ws.Range("B9").Value = 0.1132075471698113
$ws.Range("D9").Value = 0.01886792452830189
$ws.Range("F9").Value = 0.0830188679245283
$ws.Range("J9").Value = 0.07924528301886792
$ws.Range("O9").Value = 0.02264150943396226
$ws.Range("Q9").Value = 0.2
$ws.Range("R9").Value = 0.1094339622641509
$ws.Range("S9").Value = 0.3735849056603773
$ws.Range("B10").Value = 0.1181948424068768
$ws.Range("D10").Value = 0.02363896848137536
$ws.Range("E10").Value = 0.001432664756446991
$ws.Range("F10").Value = 0.07521489971346705
$ws.Range("J10").Value = 0.0988538681948424
$ws.Range("O10").Value = 0.02722063037249284
$ws.Range("Q10").Value = 0.2206303724928367
$ws.Range("R10").Value = 0.1010028653295129
$ws.Range("S10").Value = 0.333810888252149
$ws.Range("G11").Value = 0.1437308868501529
$ws.Range("J11").Value = 0.07339449541284404
$ws.Range("K11").Value = 0.2171253822629969
$ws.Range("L11").Value = 0.5565749235474006
$ws.Range("S11").Value = 0.009174311926605505
$ws.Range("G12").Value = 0.7448979591836735
$ws.Range("J12").Value = 0.1887755102040816
$ws.Range("K12").Value = 0.00510204081632653
$ws.Range("L12").Value = 0.05102040816326531
$ws.Range("S12").Value = 0.01020408163265306
$ws.Range("G13").Value = 0.7441860465116279
$ws.Range("J13").Value = 0.186046511627907
$ws.Range("S13").Value = 0.06976744186046512
$ws.Range("F15").Value = 0.0231023102310231
$ws.Range("H15").Value = 0.1254125412541254
$ws.Range("I15").Value = 0.04950495049504951
$ws.Range("J15").Value = 0.3564356435643564
$ws.Range("K15").Value = 0.08580858085808581
$ws.Range("O15").Value = 0.066006600660066
$ws.Range("S15").Value = 0.2937293729372937
$ws.Range("F16").Value = 0.02008032128514056
$ws.Range("H16").Value = 0.1726907630522088
$ws.Range("I16").Value = 0.07228915662650602
$ws.Range("J16").Value = 0.429718875502008
$ws.Range("K16").Value = 0.1004016064257028
$ws.Range("M16").Value = 0.02008032128514056
$ws.Range("O16").Value = 0.06024096385542169
$ws.Range("S16").Value = 0.1244979919678715
$ws.Range("F17").Value = 0.02033271719038817
$ws.Range("H17").Value = 0.1645101663585952
$ws.Range("I17").Value = 0.1386321626617375
$ws.Range("J17").Value = 0.4121996303142329
$ws.Range("K17").Value = 0.08872458410351201
$ws.Range("M17").Value = 0.01478743068391867
$ws.Range("N17").Value = 0.001848428835489834
$ws.Range("O17").Value = 0.05545286506469501
$ws.Range("S17").Value = 0.1035120147874307
$ws.Range("F18").Value = 0.02702702702702703
$ws.Range("H18").Value = 0.1583011583011583
$ws.Range("I18").Value = 0.08880308880308881
$ws.Range("J18").Value = 0.4015444015444015
$ws.Range("K18").Value = 0.08494208494208494
$ws.Range("M18").Value = 0.0193050193050193
$ws.Range("O18").Value = 0.07335907335907337
$ws.Range("S18").Value = 0.1467181467181467
$ws.Range("F19").Value = 0.02083333333333333
$ws.Range("H19").Value = 0.2053571428571428
$ws.Range("I19").Value = 0.09821428571428571
$ws.Range("J19").Value = 0.3690476190476191
$ws.Range("K19").Value = 0.0974702380952381
$ws.Range("M19").Value = 0.0193452380952381
$ws.Range("N19").Value = 0.000744047619047619
$ws.Range("O19").Value = 0.08556547619047619
$ws.Range("S19").Value = 0.103422619047619
